# The template field call "{m: name.asImage().setWidth(200)}" currently
# ends with a single run whose text is ")}" (the closing parenthesis of
# setWidth(...) immediately followed by the closing brace of the field).
#
# The migrated parser (TokenIteratorFieldRewriterSplit) now expects the
# closing parenthesis and the closing brace to live in two separate runs,
# so split the ")}" run into a ")" run and a "}" run while keeping the
# same run formatting (w:rPr = <w:lang w:val="en-US"/>).

$d = $word.ActiveDocument

# Locate the paragraph holding the "...setWidth(200)}" field call.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*.setWidth(*)}*") {
        $target = $p
    }
}

if ($target -ne $null) {
    $searchRange = $target.Range
    $found = $searchRange.Find.Execute(")}", $true, $false, $false, $false, $false, `
                                        $true, 1, $false, "", 0)

    if ($found) {
        # Work on a fresh Range with the same bounds as the Find hit so the
        # subsequent InsertXML only replaces the matched "}" text.
        $closingRange = $d.Range($searchRange.Start, $searchRange.End)

        $newRunsXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                      '<w:r w:rsidR="00833972"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>)</w:t></w:r>' +
                      '<w:r w:rsidR="00833972"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">}</w:t></w:r>' +
                      '</w:p>'

        $closingRange.InsertXML($newRunsXml)
    }
}
